# Update cryptocurrency price/volume table (Coin, Link, Price, Volume(1h))
# with the latest scraped values. A couple of rows (Litecoin/NEARProtocol and
# ARBITRUM/BabyDogeCoin) also swapped rank order, so their Coin/Link/Price/Volume
# cells are rewritten in place rather than just the numeric columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.241.54"
$ws.Range("E2").Value = "  +0.51%  "
# Row 3
$ws.Range("D3").Value = "2.479.00"
$ws.Range("E3").Value = "  +1.07%  "
# Row 4
$ws.Range("E4").Value = "  -0.02%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.13%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.97%  "
# Row 7
$ws.Range("E7").Value = "  -0.06%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.516"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.77%  "
# Row 9
$ws.Range("D9").Value = "2.478.63"
$ws.Range("E9").Value = "  +1.03%  "
# Row 10
$ws.Range("E10").Value = "  +3.15%  "
# Row 11
$ws.Range("E11").Value = "  +1.42%  "
# Row 12
$ws.Range("E12").Value = "  +1.45%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.334"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.66%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.12%  "
# Row 16
$ws.Range("D16").Value = "67.083.41"
$ws.Range("E16").Value = "  +0.35%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000171"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.66%  "
# Row 18
$ws.Range("D18").Value = "2.490.95"
$ws.Range("E18").Value = "  -0.09%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.43%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.56%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.86%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.44%  "
# Row 23
$ws.Range("E23").Value = "  +0.04%  "
# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.72%  "
# Row 25
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.03%  "
# Row 26
$ws.Range("E26").Value = "  +3.06%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.88%  "
# Row 28
$ws.Range("D28").Value = "2.574.60"
$ws.Range("E28").Value = "  +0.10%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.33%  "
# Row 30
$ws.Range("D30").Value = "0.0₃0909"
$ws.Range("E30").Value = "  +1.41%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "512.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.51%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.70%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.71%  "
# Row 34
$ws.Range("E34").Value = "  +0.00%  "
# Row 35
$ws.Range("E35").Value = "  -0.03%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.16%  "
# Row 37
$ws.Range("E37").Value = "  +1.38%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.76%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.18"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.36%  "
# Row 40
$ws.Range("E40").Value = "  -0.10%  "
# Row 41
$ws.Range("E41").Value = "  -0.01%  "
# Row 42
$ws.Range("E42").Value = "  +2.18%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.330"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.22%  "
# Row 44
$ws.Range("E44").Value = "  +2.06%  "
# Row 45
$ws.Range("E45").Value = "  +3.41%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.44%  "
# Row 47
$ws.Range("E47").Value = "  +0.44%  "
# Row 48
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0258"
$ws.Range("E48").Value = "  +3.45%  "
# Row 49
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.516"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.61%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0735"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.73%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.29%  "
